$d = $word.ActiveDocument

$replacements = @(
    @("23×87=2001", "52×40=2080"),
    @("57×24=1368", "26×83=2158"),
    @("36×37=1332", "11×31=341"),
    @("94×86=8084", "95×22=2090"),
    @("48×47=2256", "77×67=5159"),
    @("22×12=264",  "85×12=1020"),
    @("97×66=6402", "98×24=2352"),
    @("16×32=512",  "37×99=3663"),
    @("83×73=6059", "75×14=1050"),
    @("62×64=3968", "62×76=4712"),
    @("18×68=1224", "85×16=1360"),
    @("69×30=2070", "54×81=4374"),
    @("98×11=1078", "12×70=840"),
    @("22×80=1760", "55×65=3575"),
    @("79×19=1501", "35×46=1610"),
    @("29×48=1392", "12×54=648"),
    @("23×33=759",  "47×23=1081"),
    @("33×33=1089", "57×46=2622"),
    @("28×57=1596", "38×27=1026"),
    @("63×74=4662", "52×31=1612"),
    @("74×37=2738", "52×17=884"),
    @("14×47=658",  "17×43=731"),
    @("41×60=2460", "57×98=5586"),
    @("86×50=4300", "16×71=1136"),
    @("93×46=4278", "57×73=4161")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
